$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("daily")

$ws.Range("A5").Value = "AUD_USD"
$ws.Range("A6").Value = "EUR_JPY"
$ws.Range("A7").Value = "GBP_JPY"
$ws.Range("A8").Value = "NZD_USD"

$ws.Range("A9").Select()
